# Add measurement information (documentation): record the interpreter
# version used to produce the numbers on each benchmark sheet, and add a
# new "configuration" sheet describing the machine the measurements were
# taken on.

$wb = $excel.ActiveWorkbook

$xlRight = -4152
$xlLeft  = -4131

# --- guile sheet: record the Guile version used ----------------------------
$guile = $wb.Worksheets.Item("guile")
$guile.Range("A28").Value = "Version:"
$guile.Range("A28").Font.Bold = $true
$guile.Range("A28").HorizontalAlignment = $xlRight
$guile.Range("B28").Value = "GNU Guile 2.0.11"
$guile.Range("B28").HorizontalAlignment = $xlLeft

# --- larceny sheet: record the Larceny version used ------------------------
$larceny = $wb.Worksheets.Item("larceny")
$larceny.Range("A28").Value = "Version:"
$larceny.Range("A28").Font.Bold = $true
$larceny.Range("A28").HorizontalAlignment = $xlRight
$larceny.Range("B28").Value = "0.98 ""General Ripper"""
$larceny.Range("B28").HorizontalAlignment = $xlLeft

# --- new "configuration" sheet with the test machine's specs --------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$config = $wb.Worksheets.Add($null, $lastSheet)
$config.Name = "configuration"
$config.Columns.Item(2).ColumnWidth = 24.14

$config.Range("B2").Value = "MacBook Air (Mid 2011)"
$config.Range("B3").Value = "1.7 GHz Intel Core i5"
$config.Range("B4").Value = "4 GB @ 1333 MHz DDR3"

$config.Range("A3").Value = "CPU:"
$config.Range("A4").Value = "RAM:"
$config.Range("A5").Value = "OS:"
$config.Range("A2").Value = "Notebook:"

$config.Range("B5").Value = "Mac OS 10.10.4"

$config.Range("A2:A5").Font.Bold = $true
$config.Range("A2:A5").HorizontalAlignment = $xlRight

$config.Activate()
